$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 3 new product rows (15,16,17) before the old "totals" row (old
#    row 15, new row 18), copying the formatting of row 14 (last product
#    row) so the new rows look identical to the existing ones.
# ---------------------------------------------------------------------------
$ws.Rows("15:17").Insert()

$ws.Range("A14:Q14").Copy()
$ws.Range("A15:Q17").PasteSpecial(-4122)   # xlPasteFormats

$ws.Rows(15).RowHeight = 24.75
$ws.Rows(16).RowHeight = 25.5
$ws.Rows(17).RowHeight = 25.5

$ws.Range("A15:B15").Merge()
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

$ws.Range("A16:B16").Merge()
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()

$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

# Fix the style of the "transactions" column (Q) on the new rows - the
# format-only paste leaves it one style short of row 14's Q-column style.
$ws.Range("Q15").Style = $ws.Range("Q14").Style
$ws.Range("Q16").Style = $ws.Range("Q14").Style
$ws.Range("Q17").Style = $ws.Range("Q14").Style

# ---------------------------------------------------------------------------
# 2. Re-write the whole product table (rows 7-17) with the new, reordered
#    (alphabetical) list of 11 items - 3 brand-new items plus the 8
#    pre-existing ones.
# ---------------------------------------------------------------------------
$items = @(
    @{ Row=7;  Num=1;  Name="ABIMOL 300MG 5 RECTAL SUPP.";        Stock="8:0"; Reorder=1; Price="15.00";  Sell="15.0000";  Trans="1:0" },
    @{ Row=8;  Num=2;  Name="BI-KETOGESIC 150 MG 30 TAB.";        Stock="1:0"; Reorder=1; Price="81.00";  Sell="26.7300";  Trans="0:1" },
    @{ Row=9;  Num=3;  Name="BI-PROFENID 150MG 20 SCORED TABS."; Stock="0:1"; Reorder=1; Price="54.00";  Sell="54.0000";  Trans="1:0" },
    @{ Row=10; Num=4;  Name="CETAFEN PLUS 30 TAB.";               Stock="1:1"; Reorder=1; Price="81.00";  Sell="26.7300";  Trans="0:1" },
    @{ Row=11; Num=5;  Name="CLAVIMOX 1 GM 12 F.C.TABS.";         Stock="1:2"; Reorder=1; Price="130.00"; Sell="42.9000";  Trans="0:1" },
    @{ Row=12; Num=6;  Name="HELI-CURE 14 ENTERIC COATED TAB";    Stock="2:0"; Reorder=1; Price="240.00"; Sell="240.0000"; Trans="1:0" },
    @{ Row=13; Num=7;  Name="HIBIOTIC 1GM 16 TAB";                 Stock="1:1"; Reorder=1; Price="173.00"; Sell="86.5000";  Trans="0:1" },
    @{ Row=14; Num=8;  Name="MOBITIL 15MG/1.5ML 3 AMP.";          Stock="1:1"; Reorder=1; Price="39.00";  Sell="12.8700";  Trans="0:1" },
    @{ Row=15; Num=9;  Name="ZITHOTRAC 500 MG 3 TAB";              Stock="0:0"; Reorder=0; Price="50.00";  Sell="50.0000";  Trans="1:0" },
    @{ Row=16; Num=10; Name="سرنجات 3 سم";                         Stock="0:0"; Reorder=0; Price="2.00";   Sell="2.0000";   Trans="1:0" },
    @{ Row=17; Num=11; Name="سرنجات 5 سم";                         Stock="0:0"; Reorder=0; Price="3.00";   Sell="3.0000";   Trans="1:0" }
)

foreach ($it in $items) {
    $r = $it.Row
    $ws.Range("A$r").Value = $it.Num
    $ws.Range("C$r").Value = $it.Name
    $ws.Range("H$r").Value = $it.Stock
    $ws.Range("L$r").Value = $it.Reorder
    $ws.Range("N$r").Value = $it.Price
    $ws.Range("P$r").Value = $it.Sell
    $ws.Range("Q$r").Value = $it.Trans
}

# ---------------------------------------------------------------------------
# 3. Update the total (was row 15, now row 18) and the timestamp in the
#    footer (was row 16, now row 19).
# ---------------------------------------------------------------------------
$ws.Range("P18").Value = 559.73

$ws.Range("A19").Value = "Sunday, 25 May, 2025 12:15 PM"

Write-Output "done"
